$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.119630336761475
$ws.Range("B1").Value = 1.412433385848999
$ws.Range("C1").Value = 1.377718448638916
$ws.Range("D1").Value = 1.578975081443787
$ws.Range("E1").Value = 1.296739816665649
